$d = $word.ActiveDocument

# 1. Collapse "## Variables: 11" + linebreak + the "Nombre_Empresa" debug line into
#    a single "## Variables: 10" line (the Nombre_Empresa variable/column was removed
#    from the dataset, so its glimpse() output line disappears too).
$searchText = "## Variables: 11`v## `$ Nombre_Empresa                <chr> `"Bentancur Costabarria Alvaro`", ..."
$replaceText = "## Variables: 10"
$d.Content.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# 2. Update the column count bullet.
$d.Content.Find.Execute("11 columnas", $false, $false, $false, $false, $false, $true, 1, $false, "10 columnas", 2) | Out-Null

# 3. Remove "Nombre_Empresa, " from the variable-name listing bullet.
$d.Content.Find.Execute("Nombre variables: Nombre_Empresa, Localidad", $false, $false, $false, $false, $false, $true, 1, $false, "Nombre variables: Localidad", 2) | Out-Null
